# [이용섭] Add - [Prop] Palette 1종 추가
# Fixed - [Palette/Prop] 프롭 테이블화 및 Storage 존타입 구분화.
#
# Rename the storage-fixture rows from internal asset codes (SM_FUR_*) to
# human-readable zone names, add the new "Palette" prop row data, and bump
# the quantity/id values on the Palette / Shelf Stand rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the ZoneType column (D5:D15) with the new descriptive names.
$ws.Range("D5").Value  = "cashier's desk 1"
$ws.Range("D6").Value  = "cashier's desk 2"
$ws.Range("D7").Value  = "Normal Display 1"
$ws.Range("D8").Value  = "refrigeration 1"
$ws.Range("D9").Value  = "refrigeration 2"
$ws.Range("D10").Value = "Frozen Display 1"
$ws.Range("D11").Value = "Normal Display 2"
$ws.Range("D12").Value = "Fresh Display 1"
$ws.Range("D13").Value = "Fresh Display 1"
$ws.Range("D14").Value = "Palette"
$ws.Range("D15").Value = "Shelf Stand 1"

# Update quantity / id values for the Palette (row 14) and Shelf Stand (row 15) rows.
$ws.Range("B14").Value = 3
$ws.Range("H14").Value = 4012
$ws.Range("I14").Value = 1011
$ws.Range("B15").Value = 2

# Restore the selection to the data range as saved by the author.
$ws.Range("A5:I15").Select()
